$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats - used below to copy just the cell formatting (number
# format / font / border / alignment) without disturbing values.
$xlPasteFormats = -4122

# Append the new row (18) of raw/clean SSA data for June 19th, 2020.
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 44001
$ws.Range("C20").Value = 170485
$ws.Range("D20").Value = 233137
$ws.Range("E20").Value = 62245
$ws.Range("F20").Value = 20394
$ws.Range("G20").Value = 31.72

# Row 20 is now the last row in the log, so its date cell takes on the
# distinctive "final row" date-only formatting that row 19 used to have.
$ws.Range("B19").Copy()
$ws.Range("B20").PasteSpecial($xlPasteFormats)

# Row 19 is no longer the last row, so its date cell reverts to the
# standard date/time formatting shared by the rest of column B.
$ws.Range("B18").Copy()
$ws.Range("B19").PasteSpecial($xlPasteFormats)

# Give the new row's index cell (column A) the same formatting used by
# the rest of column A (bold, centered, bordered).
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false
